# Add 2022-Q3 data
# 1) Insert a new quarterly sheet "2022-Q3" right after the "总计" (summary) sheet,
#    cloned from "2022-Q2" so it keeps the same look & feel / styling.
# 2) Populate it with the new fund-holding rows for 2022-Q3.
# 3) Update the "总计" summary sheet with a new top row for 2022-Q3 and shift
#    the rest of the quarters down by one row.
# 4) Keep the last sheet ("2020-Q4") as the active tab, matching the original file.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Step 1: create the "2022-Q3" worksheet right after "总计" by cloning
# "2022-Q2" (same header wording / cell styles as the other quarter sheets).
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item("总计")
$template = $wb.Worksheets.Item("2022-Q2")
$template.Copy($null, $total)

$q3 = $wb.Worksheets.Item(2)
$q3.Name = "2022-Q3"

# The template sheet only has 2 data rows (rows 2-3); we need 7 data rows
# (rows 2-8), so extend the formatting of row 2 down through row 8.
$q3.Range("A2:H2").Copy()
$q3.Range("A3:H8").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# Step 2: fill in the 2022-Q3 fund holdings table.
# ---------------------------------------------------------------------------
$q3.Range("B1").Value = "基金代码"
$q3.Range("C1").Value = "基金名称"
$q3.Range("D1").Value = "基金规模"
$q3.Range("E1").Value = "股票总仓位"
$q3.Range("F1").Value = "仓位占比"
$q3.Range("G1").Value = "持有市值(亿元)"
$q3.Range("H1").Value = "仓位排名"

$q3.Range("A2").Value = 0
$q3.Range("B2").Value = "'011429"
$q3.Range("C2").Value = "前海开源民裕进取混合"
$q3.Range("D2").Value = "'2.53"
$q3.Range("E2").Value = "'60.79"
$q3.Range("F2").Value = "'2.99"
$q3.Range("G2").Value = "'0.0756"
$q3.Range("H2").Value = 9

$q3.Range("A3").Value = 1
$q3.Range("B3").Value = "'006195"
$q3.Range("C3").Value = "国金量化多因子股票A"
$q3.Range("D3").Value = "'5.47"
$q3.Range("E3").Value = "'84.51"
$q3.Range("F3").Value = "'0.92"
$q3.Range("G3").Value = "'0.0503"
$q3.Range("H3").Value = 3

$q3.Range("A4").Value = 2
$q3.Range("B4").Value = "'011243"
$q3.Range("C4").Value = "万家惠裕回报6个月持有期混合A"
$q3.Range("D4").Value = "'1.54"
$q3.Range("E4").Value = "'27.67"
$q3.Range("F4").Value = "'1.34"
$q3.Range("G4").Value = "'0.0206"
$q3.Range("H4").Value = 3

$q3.Range("A5").Value = 3
$q3.Range("B5").Value = "'002872"
$q3.Range("C5").Value = "华夏智胜价值成长股票C"
$q3.Range("D5").Value = "'2.13"
$q3.Range("E5").Value = "'93.39"
$q3.Range("F5").Value = "'0.80"
$q3.Range("G5").Value = "'0.0170"
$q3.Range("H5").Value = 9

$q3.Range("A6").Value = 4
$q3.Range("B6").Value = "'002871"
$q3.Range("C6").Value = "华夏智胜价值成长股票A"
$q3.Range("D6").Value = "'0.86"
$q3.Range("E6").Value = "'93.39"
$q3.Range("F6").Value = "'0.80"
$q3.Range("G6").Value = "'0.0069"
$q3.Range("H6").Value = 9

$q3.Range("A7").Value = 5
$q3.Range("B7").Value = "'011244"
$q3.Range("C7").Value = "万家惠裕回报6个月持有期混合C"
$q3.Range("D7").Value = "'0.12"
$q3.Range("E7").Value = "'27.67"
$q3.Range("F7").Value = "'1.34"
$q3.Range("G7").Value = "'0.0016"
$q3.Range("H7").Value = 3

$q3.Range("A8").Value = 6
$q3.Range("B8").Value = "'016858"
$q3.Range("C8").Value = "国金量化多因子股票C"
$q3.Range("D8").Value = "'0.00"
$q3.Range("E8").Value = "'84.51"
$q3.Range("F8").Value = "'0.92"
$q3.Range("G8").Value = 0
$q3.Range("H8").Value = 3

# ---------------------------------------------------------------------------
# Step 3: update the "总计" summary sheet - insert a new row for 2022-Q3
# above the existing data and renumber the index column.
# ---------------------------------------------------------------------------
$total.Rows(2).Insert()
$total.Range("A2:D2").Select()
$total.Range("A3:D3").Copy()
$total.Range("A2:D2").PasteSpecial(-4122)

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q3"
$total.Range("C2").Value = 7
$total.Range("D2").Value = 0.17

$total.Range("A3").Value = 1
$total.Range("B3").Value = "2022-Q2"
$total.Range("C3").Value = 2
$total.Range("D3").Value = 0.02

$total.Range("A4").Value = 2
$total.Range("B4").Value = "2022-Q1"
$total.Range("C4").Value = 4
$total.Range("D4").Value = 0.08

$total.Range("A5").Value = 3
$total.Range("B5").Value = "2021-Q4"
$total.Range("C5").Value = 3
$total.Range("D5").Value = 1.94

$total.Range("A6").Value = 4
$total.Range("B6").Value = "2021-Q3"
$total.Range("C6").Value = 1
$total.Range("D6").Value = 3.01

$total.Range("A7").Value = 5
$total.Range("B7").Value = "2021-Q2"
$total.Range("C7").Value = 1
$total.Range("D7").Value = 2.76

$total.Range("A8").Value = 6
$total.Range("B8").Value = "2021-Q1"
$total.Range("C8").Value = 1
$total.Range("D8").Value = 2.62

$total.Range("A9").Value = 7
$total.Range("B9").Value = "2020-Q4"
$total.Range("C9").Value = 2
$total.Range("D9").Value = 0.02

# ---------------------------------------------------------------------------
# Step 4: keep "2020-Q4" as the active/selected sheet, like the original file.
# ---------------------------------------------------------------------------
$last = $wb.Worksheets.Item($wb.Worksheets.Count)
$last.Activate()
$last.Range("A1").Select()
